$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Logic change for "Logged in User": the Runmode column (E) should no longer
# ever be skipped ("No") - every testcase now runs, so flip every "No" to "Yes".
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Range("E$r")
    if ($cell.Value2 -eq "No") {
        $cell.Value = "Yes"
    }
}

# Move the active selection up to E19 (was E22).
$ws.Range("E19").Select()
